$d = $word.ActiveDocument

$d.Content.Find.Execute('${{unit_price}}', $true, $false, $false, $false, $false,
                         $true, 1, $false, '$ {{unit_price}}', 2)

$d.Content.Find.Execute('For longer probes please add ${{length_adder}}', $true, $false, $false, $false, $false,
                         $true, 1, $false, 'For longer probes please add $ {{length_adder}}', 2)
